# Applies the stock-report correction described by the commit diff.
# Each entry below is a cell reference -> the corrected value that should
# be written into it (quantities, recalculated values, subtotal/grand
# total rollups, and a couple of swapped rows where the two consignment
# codes traded places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    "F13" = 11;
    "G13" = 18496.72;
    "F14" = 8;
    "G14" = 11074.48;
    "B15" = 75273.67999999999;
    "F24" = 144;
    "G24" = 16598.88;
    "F26" = 318;
    "G26" = 16291.14;
    "F27" = 75;
    "G27" = 1920.75;
    "F33" = 23;
    "G33" = 824.78;
    "F35" = 60;
    "G35" = 3073.8;
    "B40" = 52505.23;
    "F42" = 65;
    "G42" = 12790.05;
    "F47" = 181;
    "G47" = 34913.09;
    "F50" = 21;
    "G50" = 741.72;
    "F51" = 77;
    "G51" = 1750.98;
    "F54" = 212;
    "G54" = 11893.2;
    "F57" = 121;
    "G57" = 11318.34;
    "F59" = 23;
    "G59" = 377.43;
    "F61" = 102;
    "G61" = 5687.52;
    "F64" = 45;
    "G64" = 3506.85;
    "F67" = 189;
    "G67" = 49277.97;
    "B72" = 171080.57;
    "F112" = 180;
    "G112" = 20192.4;
    "B129" = 67164.91;
    "B132" = 65258;
    "B133" = 64196;
    "F181" = 22;
    "G181" = 6380.22;
    "F190" = 42;
    "G190" = 3444;
    "B199" = 54978.96;
    "F212" = 47;
    "G212" = 3045.6;
    "B214" = 3045.6;
    "F226" = 29;
    "G226" = 1832.22;
    "F228" = 278;
    "G228" = 5143;
    "F233" = 21;
    "G233" = 2406.6;
    "B235" = 12161.35;
    "F269" = 0;
    "G269" = 0;
    "F270" = 51;
    "G270" = 1776.84;
    "B301" = 95180.59;
    "B310" = 63520;
    "E310" = 153.4;
    "F310" = 35;
    "G310" = 5049.8;
    "B311" = 55373;
    "E311" = 163.62;
    "F311" = -94;
    "G311" = -13562.32;
    "B312" = 63531;
    "E312" = 152.53;
    "F312" = 24;
    "G312" = 3443.52;
    "B313" = 57802;
    "E313" = 162.71;
    "F313" = -79;
    "G313" = -11334.92;
    "F355" = 122;
    "G355" = 9104.860000000001;
    "F357" = 186;
    "G357" = 26903.04;
    "B362" = 70546.91;
    "F367" = 189;
    "G367" = 26571.51;
    "B369" = 58497.89;
    "F376" = 159;
    "G376" = 26392.41;
    "B378" = 47380.16;
    "F393" = 348;
    "G393" = 33616.8;
    "B395" = 49677.12;
    "F402" = 112;
    "G402" = 2853.76;
    "F412" = 72;
    "G412" = 7406.64;
    "F419" = 64;
    "G419" = 3684.48;
    "F421" = 50;
    "G421" = 2735;
    "B423" = 155478.04;
    "F429" = 50;
    "G429" = 4179.5;
    "B433" = 22479.62;
    "F436" = 199;
    "G436" = 9209.719999999999;
    "F439" = 80;
    "G439" = 771.2;
    "B444" = 20305.5;
    "F462" = 141;
    "G462" = 38051.67;
    "B464" = 81724.59;
    "B485" = 64810;
    "E485" = 291.22;
    "F485" = 0;
    "G485" = 0;
    "B486" = 53319;
    "E486" = 310.64;
    "F486" = -6;
    "G486" = -1643.52;
    "F517" = 170;
    "G517" = 16977.9;
    "F525" = 345;
    "G525" = 18933.6;
    "F530" = 15;
    "G530" = 1325.1;
    "B531" = 106924.65;
    "F564" = 130;
    "G564" = 15840.5;
    "F566" = 7;
    "G566" = 562.66;
    "B567" = 17754.38;
    "F570" = 4;
    "G570" = 770.24;
    "B579" = 12792.54;
    "F621" = 92;
    "G621" = 14229.64;
    "F623" = 0;
    "G623" = 0;
    "F628" = 453;
    "G628" = 46618.23;
    "B634" = 190614.96;
    "F680" = 352;
    "G680" = 57414.72;
    "B686" = 58427.27;
    "F695" = 25;
    "G695" = 3133.75;
    "B697" = 9273.389999999999;
    "F712" = 104;
    "G712" = 4075.76;
    "F717" = 9;
    "G717" = 4804.65;
    "B719" = 55181.44;
    "B724" = 2184354.55;
    "B725" = 2184354.55
}

foreach ($cellRef in $cellUpdates.Keys) {
    $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
}

Write-Output ("Updated {0} cells." -f $cellUpdates.Count)
